$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Header cells (text dates) ---
function Set-TextValue($ws, $addr, $val) {
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $val
    $ws.Range($addr).ClearFormats()
}

function Set-BlankCell($ws, $addr) {
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).ClearFormats()
}

Set-TextValue $ws1 "FK4" "10/01/2025"
Set-TextValue $ws1 "FL4" "11/01/2025"
Set-TextValue $ws2 "EY4" "10/01/2025"
Set-TextValue $ws2 "EZ4" "11/01/2025"

# --- Sheet1 (TABLE_1) data rows ---
$ws1.Range("FK5").Value = 11201.1
$ws1.Range("FL5").Value = 11312.1
$ws1.Range("FK6").Value = 182.5
$ws1.Range("FL6").Value = 183.7
$ws1.Range("FK7").Value = 28.1
$ws1.Range("FL7").Value = 28.4
$ws1.Range("FK8").Value = 201.3
$ws1.Range("FL8").Value = 204.1
$ws1.Range("FK9").Value = 103.7
$ws1.Range("FL9").Value = 103.9
$ws1.Range("FK10").Value = 1312.7
$ws1.Range("FL10").Value = 1334.7
$ws1.Range("FK11").Value = 246.1
$ws1.Range("FL11").Value = 248.7
$ws1.Range("FK12").Value = 127.3
$ws1.Range("FL12").Value = 128
$ws1.Range("FK13").Value = 39
$ws1.Range("FL13").Value = 39.4
Set-BlankCell $ws1 "FK14"
Set-BlankCell $ws1 "FL14"
$ws1.Range("FK15").Value = 506.3
$ws1.Range("FL15").Value = 505.8
$ws1.Range("FK16").Value = 354.6
$ws1.Range("FL16").Value = 355.4
Set-BlankCell $ws1 "FK17"
Set-BlankCell $ws1 "FL17"
$ws1.Range("FK18").Value = 62.9
$ws1.Range("FL18").Value = 63.5
$ws1.Range("FK19").Value = 471.8
$ws1.Range("FL19").Value = 473.6
$ws1.Range("FK20").Value = 235.6
$ws1.Range("FL20").Value = 237.8
$ws1.Range("FK21").Value = 152.4
$ws1.Range("FL21").Value = 153.6
$ws1.Range("FK22").Value = 137.6
$ws1.Range("FL22").Value = 138.7
$ws1.Range("FK23").Value = 167.1
$ws1.Range("FL23").Value = 168.5
$ws1.Range("FK24").Value = 154.1
$ws1.Range("FL24").Value = 155.3
$ws1.Range("FK25").Value = 49.1
$ws1.Range("FL25").Value = 49.5
$ws1.Range("FK26").Value = 237.4
$ws1.Range("FL26").Value = 240.7
$ws1.Range("FK27").Value = 245.5
$ws1.Range("FL27").Value = 248.1
$ws1.Range("FK28").Value = 330.2
$ws1.Range("FL28").Value = 334.2
$ws1.Range("FK29").Value = 211.3
$ws1.Range("FL29").Value = 213.7
$ws1.Range("FK30").Value = 105.4
$ws1.Range("FL30").Value = 105.9
Set-BlankCell $ws1 "FK31"
Set-BlankCell $ws1 "FL31"
$ws1.Range("FK32").Value = 41.4
$ws1.Range("FL32").Value = 41.9
$ws1.Range("FK33").Value = 98.7
$ws1.Range("FL33").Value = 99.2
$ws1.Range("FK34").Value = 89.5
$ws1.Range("FL34").Value = 90.9
$ws1.Range("FK35").Value = 45.7
$ws1.Range("FL35").Value = 46
$ws1.Range("FK36").Value = 336.8
$ws1.Range("FL36").Value = 339.8
$ws1.Range("FK37").Value = 84.5
$ws1.Range("FL37").Value = 86
$ws1.Range("FK38").Value = 671.3
$ws1.Range("FL38").Value = 676.4
$ws1.Range("FK39").Value = 350.1
$ws1.Range("FL39").Value = 352.1
$ws1.Range("FK40").Value = 40.4
$ws1.Range("FL40").Value = 40.9
$ws1.Range("FK41").Value = 392.5
$ws1.Range("FL41").Value = 395.1
$ws1.Range("FK42").Value = 158.3
$ws1.Range("FL42").Value = 159.4
$ws1.Range("FK43").Value = 141.4
$ws1.Range("FL43").Value = 142.5
$ws1.Range("FK44").Value = 343.1
$ws1.Range("FL44").Value = 348.4
$ws1.Range("FK45").Value = 29.5
$ws1.Range("FL45").Value = 30
$ws1.Range("FK46").Value = 168.2
$ws1.Range("FL46").Value = 169.3
$ws1.Range("FK47").Value = 39.1
$ws1.Range("FL47").Value = 39.4
$ws1.Range("FK48").Value = 222.5
$ws1.Range("FL48").Value = 224.6
$ws1.Range("FK49").Value = 1226.4
$ws1.Range("FL49").Value = 1233.3
$ws1.Range("FK50").Value = 140.8
$ws1.Range("FL50").Value = 142
$ws1.Range("FK51").Value = 33.9
$ws1.Range("FL51").Value = 34.4
$ws1.Range("FK52").Value = 348.1
$ws1.Range("FL52").Value = 354.1
$ws1.Range("FK53").Value = 258.6
$ws1.Range("FL53").Value = 263.4
$ws1.Range("FK54").Value = 62.5
$ws1.Range("FL54").Value = 63.5
$ws1.Range("FK55").Value = 208.3
$ws1.Range("FL55").Value = 210.7
$ws1.Range("FK56").Value = 31
$ws1.Range("FL56").Value = 31

# --- Sheet2 (TABLE_2) data rows ---
$ws2.Range("EY5").Value = 0.298178692311802
$ws2.Range("EZ5").Value = 0.274795897563197
$ws2.Range("EY6").Value = 1.10803324099723
$ws2.Range("EZ6").Value = 1.10071546505228
$ws2.Range("EY7").Value = -1.40350877192982
$ws2.Range("EZ7").Value = -1.04529616724739
$ws2.Range("EY8").Value = -1.12966601178781
$ws2.Range("EZ8").Value = -0.826044703595732
$ws2.Range("EY9").Value = 0.67961165048544
$ws2.Range("EZ9").Value = 0.678294573643414
$ws2.Range("EY10").Value = 1.65724463718732
$ws2.Range("EZ10").Value = 1.83108262760355
$ws2.Range("EY11").Value = 3.35993280134396
$ws2.Range("EZ11").Value = 3.15221899626711
$ws2.Range("EY12").Value = 1.59616919393457
$ws2.Range("EZ12").Value = 1.50674068199842
$ws2.Range("EY13").Value = -0.51020408163266
$ws2.Range("EZ13").Value = -0.755667506297222
Set-BlankCell $ws2 "EY14"
Set-BlankCell $ws2 "EZ14"
$ws2.Range("EY15").Value = 1.13863363963244
$ws2.Range("EZ15").Value = 0.377058940265921
$ws2.Range("EY16").Value = -0.365271143579644
$ws2.Range("EZ16").Value = -0.420285794340151
Set-BlankCell $ws2 "EY17"
Set-BlankCell $ws2 "EZ17"
$ws2.Range("EY18").Value = -2.48062015503875
$ws2.Range("EZ18").Value = -2.00617283950617
$ws2.Range("EY19").Value = 1.96671709531012
$ws2.Range("EZ19").Value = 1.8494623655914
$ws2.Range("EY20").Value = -0.211774671749247
$ws2.Range("EZ20").Value = -1.4504765851637
$ws2.Range("EY21").Value = 1.8716577540107
$ws2.Range("EZ21").Value = 1.72185430463576
$ws2.Range("EY22").Value = -0.649819494584842
$ws2.Range("EZ22").Value = -0.64469914040115
$ws2.Range("EY23").Value = 1.45719489981785
$ws2.Range("EZ23").Value = 0.958657878969439
$ws2.Range("EY24").Value = 0.195058517555255
$ws2.Range("EZ24").Value = 0.258231116849584
$ws2.Range("EY25").Value = -0.607287449392735
$ws2.Range("EZ25").Value = -0.402414486921535
$ws2.Range("EY26").Value = 3.12771503040834
$ws2.Range("EZ26").Value = 2.90722530996153
$ws2.Range("EY27").Value = -0.324807145757211
$ws2.Range("EZ27").Value = 0.0806776926179865
$ws2.Range("EY28").Value = 2.76999688764394
$ws2.Range("EZ28").Value = 3.11632212280159
$ws2.Range("EY29").Value = -0.471031559114461
$ws2.Range("EZ29").Value = 0.0468164794007464
$ws2.Range("EY30").Value = 0.380952380952386
$ws2.Range("EZ30").Value = 0.665399239543742
Set-BlankCell $ws2 "EY31"
Set-BlankCell $ws2 "EZ31"
$ws2.Range("EY32").Value = -1.42857142857143
$ws2.Range("EZ32").Value = -1.87353629976582
$ws2.Range("EY33").Value = 0.817160367722163
$ws2.Range("EZ33").Value = 1.01832993890019
$ws2.Range("EY34").Value = 1.01580135440181
$ws2.Range("EZ34").Value = 0.887902330743631
$ws2.Range("EY35").Value = 0.883002207505531
$ws2.Range("EZ35").Value = 1.54525386313466
$ws2.Range("EY36").Value = 0.687593423019418
$ws2.Range("EZ36").Value = 0.324771183938579
$ws2.Range("EY37").Value = 6.02258469259724
$ws2.Range("EZ37").Value = 6.56753407682775
$ws2.Range("EY38").Value = 0.208986415882964
$ws2.Range("EZ38").Value = 0.266824785057822
$ws2.Range("EY39").Value = 0.806219406852868
$ws2.Range("EZ39").Value = 0.427837992013691
$ws2.Range("EY40").Value = 0.999999999999996
$ws2.Range("EZ40").Value = 1.48883374689828
$ws2.Range("EY41").Value = 0.744353182751549
$ws2.Range("EZ41").Value = 0.0506457330969981
$ws2.Range("EY42").Value = 2.2609819121447
$ws2.Range("EZ42").Value = 2.44215938303343
$ws2.Range("EY43").Value = -0.632466619817291
$ws2.Range("EZ43").Value = -1.58839779005526
$ws2.Range("EY44").Value = 1.50887573964496
$ws2.Range("EZ44").Value = 1.2496367335077
$ws2.Range("EY45").Value = -1.66666666666667
$ws2.Range("EZ45").Value = -0.662251655629137
$ws2.Range("EY46").Value = 0.83932853717025
$ws2.Range("EZ46").Value = 0.355660936573813
$ws2.Range("EY47").Value = -0.25510204081633
$ws2.Range("EZ47").Value = -0.505050505050494
$ws2.Range("EY48").Value = 2.77136258660508
$ws2.Range("EZ48").Value = 3.07480495640203
$ws2.Range("EY49").Value = 1.34699611602348
$ws2.Range("EZ49").Value = 1.09845069267972
$ws2.Range("EY50").Value = 3.45334313005145
$ws2.Range("EZ50").Value = 3.19767441860466
$ws2.Range("EY51").Value = 1.80180180180181
$ws2.Range("EZ51").Value = 2.6865671641791
$ws2.Range("EY52").Value = 3.35510688836105
$ws2.Range("EZ52").Value = 4.05524537173083
$ws2.Range("EY53").Value = -2.11960635881906
$ws2.Range("EZ53").Value = -2.08178438661711
$ws2.Range("EY54").Value = 0.644122383252827
$ws2.Range("EZ54").Value = 1.11464968152866
$ws2.Range("EY55").Value = -2.52690687880205
$ws2.Range("EZ55").Value = -1.77156177156178
$ws2.Range("EY56").Value = -0.64102564102565
$ws2.Range("EZ56").Value = -1.27388535031847

Write-Host "edit complete"